# Generate Report for Handoff
# b.md has now been handed off for localization (new xliff files generated
# for zh-cn and de-de); the report reflects the new "Ready for handoff"
# status, the new handoff file names / timestamps, and the handback
# version-mismatch error message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1363d3dd7dc2ba43ee3f7eba0a67c16ffba220ac/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d89d647a0f4056c14702ef0133ccf0f79e8c4097/e2e/b.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 06:44:38"

# ---- zh-cn sheet (row 3 = b.md) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces Excel to treat "False" as literal text instead of
# a boolean; re-applying the Normal style clears the resulting quote-prefix
# formatting flag so the cell matches the plain text cells around it.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-20 06:44:34"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---- de-de sheet (row 3 = b.md) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-20 06:44:38"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
